# "updated feedback and deliverables"
# Apply the March 2017 Deliverables Tracking updates to the
# "Hardware Development Process" sheet, plus the matching view-state
# tweaks on both sheets.

$wb  = $excel.ActiveWorkbook
$wsProf = $wb.Worksheets.Item("Professionalism")
$wsHw   = $wb.Worksheets.Item("Hardware Development Process")

# --- Data edits on "Hardware Development Process" ---

# Row 9 "Hardware Block Diagram": Due Date 3/12 -> 3/11
$wsHw.Range("D9").Value2 = 42805

# Row 10 "Product Architecture": Due Date 3/13 -> 3/11
$wsHw.Range("D10").Value2 = 42805

# Row 12 "Estimate Architecture Task Hours": add the same red-fill
# marker cell in H12 that already exists in G12.
$wsHw.Range("G12").Copy() | Out-Null
$wsHw.Range("H12").PasteSpecial(-4122) | Out-Null

# Row 13 "Estimate Task Hours": Assigned/Due Date go from the
# placeholder "ENTER DATE" text to real dates (3/11), matching the
# date format used elsewhere in the column, and the marker cell G13
# switches from the unfilled style to the red-fill style (like G12).
$wsHw.Range("D9").Copy() | Out-Null
$wsHw.Range("C13:D13").PasteSpecial(-4122) | Out-Null
$wsHw.Range("C13").Value2 = 42805
$wsHw.Range("D13").Value2 = 42805

$wsHw.Range("G12").Copy() | Out-Null
$wsHw.Range("G13").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- View-state updates ---

$wsProf.Activate()
$wsProf.Range("G6").Select() | Out-Null

$wsHw.Activate()
$wsHw.Range("E18").Select() | Out-Null
$excel.ActiveWindow.Zoom = 190
